# Adjusting NOV13/14 observation locations.
# Row 11 and Row 12 on Sheet1 swap their "Agency"/"Station ID" pairing:
#   - Row 11 becomes the USGS station, with its station id corrected to 07380260
#     (was the unused/garbled placeholder 291929089562600).
#   - Row 12 becomes the USACE station 01440 (previously sitting in row 11).
# Column C (the numeric flag) stays 0 in both rows, untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "USGS"
$ws.Range("B11").Value = "07380260"

$ws.Range("A12").Value = "USACE"
$ws.Range("B12").Value = "01440"

# Reflect the edit location as the active selection, like it would be
# right after typing the new values into A11:B12.
$ws.Range("A11:B12").Select()
